$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Build new row 28 (copy formatting from row 27, then fill in values) ---
$ws.Range("A27:K27").Copy() | Out-Null
$ws.Range("A28:K28").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Rows.Item(28).RowHeight = 120

$ws.Range("A28").Value2 = $ws.Range("A27").Value2
$ws.Range("B28").Value2 = 27
$ws.Range("C28").Value2 = $ws.Range("C27").Value2
$ws.Range("D28").Value2 = $ws.Range("D27").Value2
$ws.Range("E28").Value2 = "https://www.mitradel.gob.pa/instituciones-capacitan-para-un-regreso-seguro-al-trabajo/"
$ws.Range("F28").Value2 = "La creación de planes de mitigación de contagio por Covid-19, para un regreso seguro al trabajo, mediante la creación de los Comités de Salud en las empresas publicas y privadas, es el objetivo fundamental de las capacitaciones virtuales sobre Estrategias y Protocolos para Preservar la Higiene y Salud en el Trabajo como Prevención ante el Covid-19, dictadas por un equipo interinstitucional integrado por especialistas de los ministerios de Trabajo y Desarrollo Laboral (Mitradel), Salud y la Caja de Seguro Social."
$ws.Range("G28").Value2 = "https://www.mitradel.gob.pa"
$ws.Range("H28").Value2 = "29-06-2020"
$ws.Range("I28").Value2 = "29-06-2020"
$ws.Range("J28").Value2 = $ws.Range("J27").Value2
$ws.Range("K28").Value2 = $ws.Range("K27").Value2

# --- Hyperlinks: new row's E/G cells, plus the previously-missing G17 link ---
$ws.Hyperlinks.Add($ws.Range("E28"), "https://www.mitradel.gob.pa/instituciones-capacitan-para-un-regreso-seguro-al-trabajo/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G28"), "https://www.mitradel.gob.pa") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G17"), "https://www.mitradel.gob.pa/decretoscovid19/") | Out-Null

# Re-apply the original cell formatting so Excel's automatic "Hyperlink"
# style doesn't clobber the table's row styling on the cells touched above.
$ws.Range("E27").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("G27").Copy() | Out-Null
$ws.Range("G28").PasteSpecial(-4122) | Out-Null
$ws.Range("G4").Copy() | Out-Null
$ws.Range("G17").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# --- Grow the table (ListObject) by one row ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:K28"))

# --- Extend the category data validation down to the new row ---
$ws.Range("C2:C28").Validation.Delete()
$ws.Range("C2:C28").Validation.Add(3, 1, 1, "=Criteria")
$validation = $ws.Range("C2:C28").Validation
$validation.ErrorTitle = "Entrada no válida"
$validation.ErrorMessage = "Selecciona una categoría de la lista"
$validation.InputTitle = "Categoria"
$validation.InputMessage = "Selecciona una categoría de la lista"
$validation.ShowError = $true
$validation.ShowInput = $true

# --- Leave the view where the editor ended up working last ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("H17").Select() | Out-Null

Write-Host "done"
